# The sheet originally lists 8 MAG rows (rows 2-9). Two of those rows are
# removed entirely:
#   - row 4: even_MAG-GUT5920.fa  (s__Veillonella_A magna)
#   - the original row 6: even_MAG-GUT9303.fa  (s__Veillonella_A sp000431435)
# Deleting full rows shifts everything below each deletion up, leaving 6
# data rows (A1:G7 used range) with all remaining values unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete even_MAG-GUT5920.fa (originally row 4).
$ws.Rows(4).Delete()

# After the above shift, even_MAG-GUT9303.fa (originally row 6) is now row 5.
$ws.Rows(5).Delete()
